$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (year 2025) figures for Dados BIBI faturamento anual
$ws.Range("B9").Value = 2849787.1
$ws.Range("C9").Value = 448651.99
$ws.Range("D9").Value = 3298439.09
$ws.Range("E9").Value = 13.60194861139606
$ws.Range("F9").Value = 86.39805138860395
$ws.Range("G9").Value = -56.64002224845532
$ws.Range("H9").Value = -48.53678406262261
$ws.Range("I9").Value = 28289
$ws.Range("J9").Value = 1210
$ws.Range("K9").Value = 29499
$ws.Range("L9").Value = 20340
$ws.Range("M9").Value = 162.1651470009833
$ws.Range("N9").Value = 10.7134487934051
